$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking values (e.g. "1.006")
# are stored as text, matching the source data, not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.924.53'
$ws.Range("D3").Value = '1.643.89'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").Value = '216.05'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").Value = '0.5057'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '1.006'
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").Value = '0.2580'
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").Value = '0.06412'
$ws.Range("E9").Value = '  -0.46%  '
$ws.Range("D10").Value = '19.65'
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").Value = '0.07812'
$ws.Range("E11").Value = '  +1.27%  '
$ws.Range("D12").Value = '1.667.76'
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").Value = '4.295'
$ws.Range("E13").Value = '  +1.25%  '
$ws.Range("D14").Value = '0.5450'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").Value = '0.0₅7889'
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").Value = '65.03'
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("D17").Value = '25.986.65'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '1.006'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").Value = '199.01'
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("D20").Value = '4.412'
$ws.Range("E20").Value = '  +3.07%  '
$ws.Range("D21").Value = '9.980'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").Value = '6.014'
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("D23").Value = '1.008'
$ws.Range("D24").Value = '1.873'
$ws.Range("E24").Value = '  -3.28%  '
$ws.Range("D25").Value = '140.86'
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("D26").Value = '0.1147'
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").Value = '6.869'
$ws.Range("E27").Value = '  +2.36%  '
$ws.Range("D28").Value = '15.77'
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("D29").Value = '1.248'
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("D30").Value = '0.04986'
$ws.Range("E30").Value = '  -1.20%  '
$ws.Range("D31").Value = '3.273'
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("D32").Value = '3.204'
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("D33").Value = '1.534'
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("D34").Value = '2.375'
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").Value = '0.8962'
$ws.Range("E35").Value = '  +1.00%  '
$ws.Range("D36").Value = '2.615'
$ws.Range("E36").Value = '  -0.86%  '
$ws.Range("D37").Value = '1.145.98'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = '0.5566'
$ws.Range("E38").Value = '  -0.83%  '
$ws.Range("D39").Value = '0.01567'
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").Value = '1.009'
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("D41").Value = '5.687'
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("D42").Value = '0.8250'
$ws.Range("E42").Value = '  +2.20%  '
$ws.Range("D43").Value = '100.05'
$ws.Range("E43").Value = '  +0.35%  '
$ws.Range("D44").Value = '0.0₈120'
$ws.Range("E44").Value = '  +6.71%  '
$ws.Range("D45").Value = '1.781.35'
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").Value = '55.48'
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("D48").Value = '1.007'
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").Value = '0.05061'
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("D50").Value = '1.009'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '0.09552'
$ws.Range("E51").Value = '  +2.75%  '

# Restore the original (default) cell style now that the text values are set,
# so no residual number formatting is left on the cells.
$ws.Range("D2:D51").Style = "Normal"
